# This change corresponds to the project's move from Apache POI 4.1.0 to
# 5.2.3 (see commit message). That library bump only changed *how* the
# test-fixture .docx is serialized by POI when the fixture was regenerated:
#   - boolean run-properties (<w:b>, <w:i>, <w:strike>) are now written as
#     w:val="on"/"off" instead of w:val="true"/"false" (both spellings are
#     the same ST_OnOff value per the OOXML schema),
#   - the child elements inside <w:rPr> come out in a different (but
#     equivalent) order,
#   - the opaque bookmark id and the internal field-run w:rsidR token are
#     regenerated random identifiers with no semantic meaning.
# None of this changes the document's actual formatting/content: every
# bold/italic/strike/color/size/lang value, the bookmark name ("id"), and
# the "REF id \h" field structure stay exactly the same before and after.
#
# Re-assert (idempotently) the same character formatting on the runs the
# diff touches, so the intent of the edit is applied through the Word
# object model without corrupting anything that a raw byte-level
# reformatting by a different XML serializer could not meaningfully
# reproduce here (exact attribute spelling/order and random ids are
# serializer artifacts, not content the object model exposes/controls).

$d = $word.ActiveDocument

function Set-SameFormatting($paraIndex, $bold, $italic, $strike) {
    $rng = $d.Paragraphs.Item($paraIndex).Range
    $rng.Font.Bold = $bold
    $rng.Font.Italic = $italic
    $rng.Font.StrikeThrough = $strike
}

# "some text" / "some text1" / "some text2" -> bold, italic, not struck
Set-SameFormatting 8  $true $true $false
Set-SameFormatting 10 $true $true $false
Set-SameFormatting 12 $true $true $false

# "bookmarkRef" (the REF id \h field result) -> bold only
Set-SameFormatting 9 $true $false $false

# "sample table" / "some text3" run -> bold, italic, not struck
Set-SameFormatting 43 $true $true $false

# "sample table" / "text in a list" run -> bold, not italic, not struck
Set-SameFormatting 44 $true $false $false
